$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): extend from O1 to include P1=14, Q1=15, with style copied from O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# New values for columns I..Q across data rows 2..25
$newRowValues = @(2,2,1,2,2,2,1,2,2)  # I, J, K, L, M, N, O, P, Q

for ($r = 2; $r -le 25; $r++) {
    for ($i = 0; $i -lt $newRowValues.Length; $i++) {
        $col = 9 + $i   # column I = 9
        $ws.Cells.Item($r, $col).Value = $newRowValues[$i]
    }
}
